# Implementatieplan Vision Week 1 - add the "Evaluatie" paragraph text,
# and move the transient _GoBack bookmark from the end of the previous
# paragraph to the end of the newly written paragraph.

$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits at the end
#    of the paragraph ending in " terug." (it will be re-created at the
#    end of the new Evaluatie paragraph below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Find the (empty) paragraph right after the "Evaluatie" heading and
#    fill it in with the full write-up, preserving the exact
#    run/proofErr layout produced by Word's own spell-checker.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.Trim() -eq "Evaluatie") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}
if ($target -eq $null) {
    # Fall back to the last paragraph in the document body.
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}
$targetRange = $target.Range

$xml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>Om de functionaliteiten van de ge</w:t></w:r><w:r><w:t>&#239;</w:t></w:r><w:r><w:t>mplementeerde Image Shells te testen zullen er twee experimenten gedaan worden.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Het</w:t></w:r><w:r><w:t xml:space="preserve"> eerste</w:t></w:r><w:r><w:t xml:space="preserve"> experiment</w:t></w:r><w:r><w:t xml:space="preserve"> gaat dieper in over de keuze van het type opslag voor de pixels. Hier worden een standaard GCC Array en Vector implementatie vergeleken qua snelheid om zo te bepalen wat het beste is om te gebruiken. Verder zal hier ook de snelheid worden vergeleken met de snelheid van de basis implementatie, om te bepalen of de gekozen implementatie methode dicht bij de basis implementatie zit.</w:t></w:r><w:r><w:t xml:space="preserve"> Ten slotte wordt er ook een experiment gedaan om de correctheid van de gekozen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RGBImage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> implementatie te testen en vergelijken met de basis implementatie. Hier zal gekeken worden of de informatie van een afbeelding correct wordt opgeslagen, specifiek de afmetingen en de pixels.</w:t></w:r><w:r><w:t xml:space="preserve"> Door deze experimenten zal kunnen worden bepaalt of de door de studenten ge</w:t></w:r><w:r><w:t>&#239;</w:t></w:r><w:r><w:t xml:space="preserve">mplementeerde image </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shells</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en conversie algoritme voldoen om verdere operaties uit te voeren.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$targetIndex = $target.Index
$targetRange.InsertXML($xml)

# InsertXML inserts the new paragraph(s) *before* the original (now
# empty) paragraph mark rather than consuming it, so the original
# empty paragraph survives as a spurious extra paragraph right after
# the text we just inserted (its Range text is just the lone
# paragraph-mark character, chr(13)). Remove that leftover paragraph
# mark so the document reads exactly as intended.
$leftover = $d.Paragraphs.Item($targetIndex + 1)
if ($leftover.Range.Text -eq [string][char]13) {
    $cleanupRange = $d.Range($leftover.Range.Start - 1, $leftover.Range.End)
    $cleanupRange.Delete()
}
